# revert to v4, small visualization change
# Applies the "Тесты" sheet updates: new test-12/test-13 rows of data,
# a row-height tweak on row 3, and a small selection/viewport change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Тесты")

# --- Row 3: shrink the row height (visualization-only change) ---
$ws.Rows.Item(3).RowHeight = 60

# --- Row 13 (test #12): fill in the reverted-to-v4 augmentation test data ---
$ws.Cells.Item(13, 2).Value = 1
$ws.Cells.Item(13, 3).Value = 40
$ws.Cells.Item(13, 4).Value = 8
$ws.Cells.Item(13, 6).Value = "Небольшая аугментация датасета"
$ws.Cells.Item(13, 7).Value = "параметры теста 4, к аугментациям добавляется поворот изображения на 90, 180 или 270 градусов"
$ws.Cells.Item(13, 8).Value = "Train IoU: 0.32, Val IoU: 0.37. Нужно проверить на большем количестве эпох."
$ws.Cells.Item(13, 9).Value = "8702ea4"
$ws.Rows.Item(13).RowHeight = 45

# --- Row 14 (test #13): follow-up test data ---
$ws.Cells.Item(14, 2).Value = 1
$ws.Cells.Item(14, 3).Value = 50
$ws.Cells.Item(14, 4).Value = 8
$ws.Cells.Item(14, 7).Value = "параметры теста 12"
$ws.Cells.Item(14, 8).Value = "Train IoU: 0.30, Val IoU: 0.36. Аугментация не привела к улучшению точности."
$ws.Rows.Item(14).RowHeight = 45

# --- Update the active selection shown when the workbook is reopened ---
$ws.Activate()
$ws.Range("H15").Select()
